$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose entire contents (columns A:T) must be swapped.
$pairs = @(
    @(2, 5),
    @(3, 7),
    @(4, 6),
    @(8, 9)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("A$r1`:T$r1")
    $range2 = $ws.Range("A$r2`:T$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}
